$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 5.191972255706787
$ws.Range("B1").Value = 4.128003120422363
$ws.Range("C1").Value = 2.988499641418457
$ws.Range("D1").Value = 2.161231994628906
$ws.Range("E1").Value = 1.663712501525879
